$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.48%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'35.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'12.89%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.164"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.91%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07810"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.51%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.395"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'8.92%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'4.34%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.971"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'6.39%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9336"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.49%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'8.95%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'9.84%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08699"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.90%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03319"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'6.32%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09905"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.43%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001484"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.00%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005771"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.57%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.467"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.60%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.153"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'4.22%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'0.1327"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.61%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.323"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.80%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'5.87%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04572"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.78%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.68%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004439"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.83%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.20%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003698"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'8.88%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01777"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'13.41%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'8.30%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007766"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.90%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'6.30%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007128"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-20.62%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002187"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.06%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009181"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.91%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00005937"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.03%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'16.33%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.11%"
$ws.Range("E51").Style = "Normal"
